# Weekly update: insert 4 new price rows (new reporting date 45223) for
# "Zafiro rojo/verde" Primera/Segunda at the top of this variety block,
# pushing the existing history down by 4 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the existing row 983, shifting
# everything currently at row 983 onward down to row 987 onward.
$ws.Rows.Item(983).Resize(4).Insert()

# --- New row 983: Zafiro rojo / Primera ---
$ws.Cells.Item(983,1).Value2 = 5
$ws.Cells.Item(983,2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(983,3).Value2 = "Maule"
$ws.Cells.Item(983,4).Value2 = 45223
$ws.Cells.Item(983,5).Value2 = 7
$ws.Cells.Item(983,6).Value2 = 100112002
$ws.Cells.Item(983,7).Value2 = "Pimiento"
$ws.Cells.Item(983,8).Value2 = "Zafiro rojo"
$ws.Cells.Item(983,9).Value2 = "Primera"
$ws.Cells.Item(983,10).Value2 = 200
$ws.Cells.Item(983,11).Value2 = 35000
$ws.Cells.Item(983,12).Value2 = 35000
$ws.Cells.Item(983,13).Value2 = 35000
$ws.Cells.Item(983,14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(983,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(983,16).Value2 = 2333
$ws.Cells.Item(983,17).Value2 = 15
$ws.Cells.Item(983,18).Value2 = "Hortaliza"

# --- New row 984: Zafiro rojo / Segunda ---
$ws.Cells.Item(984,1).Value2 = 5
$ws.Cells.Item(984,2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(984,3).Value2 = "Maule"
$ws.Cells.Item(984,4).Value2 = 45223
$ws.Cells.Item(984,5).Value2 = 7
$ws.Cells.Item(984,6).Value2 = 100112002
$ws.Cells.Item(984,7).Value2 = "Pimiento"
$ws.Cells.Item(984,8).Value2 = "Zafiro rojo"
$ws.Cells.Item(984,9).Value2 = "Segunda"
$ws.Cells.Item(984,10).Value2 = 100
$ws.Cells.Item(984,11).Value2 = 30000
$ws.Cells.Item(984,12).Value2 = 30000
$ws.Cells.Item(984,13).Value2 = 30000
$ws.Cells.Item(984,14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(984,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(984,16).Value2 = 2000
$ws.Cells.Item(984,17).Value2 = 15
$ws.Cells.Item(984,18).Value2 = "Hortaliza"

# --- New row 985: Zafiro verde / Primera ---
$ws.Cells.Item(985,1).Value2 = 5
$ws.Cells.Item(985,2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(985,3).Value2 = "Maule"
$ws.Cells.Item(985,4).Value2 = 45223
$ws.Cells.Item(985,5).Value2 = 7
$ws.Cells.Item(985,6).Value2 = 100112002
$ws.Cells.Item(985,7).Value2 = "Pimiento"
$ws.Cells.Item(985,8).Value2 = "Zafiro verde"
$ws.Cells.Item(985,9).Value2 = "Primera"
$ws.Cells.Item(985,10).Value2 = 200
$ws.Cells.Item(985,11).Value2 = 30000
$ws.Cells.Item(985,12).Value2 = 30000
$ws.Cells.Item(985,13).Value2 = 30000
$ws.Cells.Item(985,14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(985,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(985,16).Value2 = 2000
$ws.Cells.Item(985,17).Value2 = 15
$ws.Cells.Item(985,18).Value2 = "Hortaliza"

# --- New row 986: Zafiro verde / Segunda ---
$ws.Cells.Item(986,1).Value2 = 5
$ws.Cells.Item(986,2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(986,3).Value2 = "Maule"
$ws.Cells.Item(986,4).Value2 = 45223
$ws.Cells.Item(986,5).Value2 = 7
$ws.Cells.Item(986,6).Value2 = 100112002
$ws.Cells.Item(986,7).Value2 = "Pimiento"
$ws.Cells.Item(986,8).Value2 = "Zafiro verde"
$ws.Cells.Item(986,9).Value2 = "Segunda"
$ws.Cells.Item(986,10).Value2 = 100
$ws.Cells.Item(986,11).Value2 = 25000
$ws.Cells.Item(986,12).Value2 = 25000
$ws.Cells.Item(986,13).Value2 = 25000
$ws.Cells.Item(986,14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(986,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(986,16).Value2 = 1667
$ws.Cells.Item(986,17).Value2 = 15
$ws.Cells.Item(986,18).Value2 = "Hortaliza"

Write-Output "Inserted 4 rows; new dimension rows: $($ws.UsedRange.Rows.Count)"
